# Apply cell value updates to the Exodus_Profits workbook per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6
$ws.Range("H6").Value = 97.75
$ws.Range("I6").Value = 97.75
$ws.Range("K6").Value = 293.25
$ws.Range("M6").Value = -181.25

# ALC row 12
$ws.Range("H12").Value = 526.4583
$ws.Range("I12").Value = 689.2727
$ws.Range("J12").Value = 388.69232
$ws.Range("K12").Value = 689.2727
$ws.Range("L12").Value = 388.69232
$ws.Range("M12").Value = -519.2727
$ws.Range("N12").Value = -728.69232

# ALC row 28
$ws.Range("H28").Value = 73418.75
$ws.Range("I28").Value = 60491.934
$ws.Range("J28").Value = 112199.2
$ws.Range("K28").Value = 60491.934
$ws.Range("L28").Value = 112199.2
$ws.Range("M28").Value = -60006.934
$ws.Range("N28").Value = -113169.2

# ALC row 80
$ws.Range("H80").Value = 5371.1
$ws.Range("I80").Value = 1433.9375
$ws.Range("K80").Value = 4301.8125
$ws.Range("M80").Value = -3303.8125

# ALC row 83
$ws.Range("H83").Value = 5371.1
$ws.Range("I83").Value = 1433.9375
$ws.Range("K83").Value = 12905.4375
$ws.Range("M83").Value = -7913.4375

# ALC row 98
$ws.Range("H98").Value = 2451.1292
$ws.Range("I98").Value = 2519.5
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 2519.5
$ws.Range("L98").Value = 400
$ws.Range("M98").Value = -1021.5
$ws.Range("N98").Value = -3396

# ALC row 107
$ws.Range("H107").Value = 938.8
$ws.Range("I107").Value = 933.3333
$ws.Range("J107").Value = 947
$ws.Range("K107").Value = 933.3333
$ws.Range("L107").Value = 947
$ws.Range("M107").Value = 986.6667
$ws.Range("N107").Value = -4787

# ALC row 122
$ws.Range("H122").Value = 2451.1292
$ws.Range("I122").Value = 2519.5
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 7558.5
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = -5108.5
$ws.Range("N122").Value = -6100

# ALC row 132
$ws.Range("H132").Value = 1609.3695
$ws.Range("I132").Value = 1419.2106
$ws.Range("J132").Value = 2512.625
$ws.Range("K132").Value = 4257.6318
$ws.Range("L132").Value = 7537.875
$ws.Range("M132").Value = -1727.6318
$ws.Range("N132").Value = -12597.875

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 2745.63
$ws.Range("I32").Value = 2406.2444
$ws.Range("J32").Value = 5800.1
$ws.Range("K32").Value = 2406.2444
$ws.Range("L32").Value = 5800.1
$ws.Range("M32").Value = -2119.2444
$ws.Range("N32").Value = -6374.1

# ARM row 61
$ws.Range("H61").Value = 1142.2609
$ws.Range("I61").Value = 739.0454999999999
$ws.Range("J61").Value = 10013
$ws.Range("K61").Value = 739.0454999999999
$ws.Range("L61").Value = 10013
$ws.Range("M61").Value = -527.0454999999999
$ws.Range("N61").Value = -10437

# ARM row 74
$ws.Range("H74").Value = 2000.5264
$ws.Range("I74").Value = 1444.2727
$ws.Range("J74").Value = 2765.375
$ws.Range("K74").Value = 1444.2727
$ws.Range("L74").Value = 2765.375
$ws.Range("M74").Value = -570.2727
$ws.Range("N74").Value = -4513.375

# ARM row 77
$ws.Range("H77").Value = 2000.5264
$ws.Range("I77").Value = 1444.2727
$ws.Range("J77").Value = 2765.375
$ws.Range("K77").Value = 7221.363499999999
$ws.Range("L77").Value = 13826.875
$ws.Range("M77").Value = -2853.363499999999
$ws.Range("N77").Value = -22562.875

# ARM row 102
$ws.Range("H102").Value = 49282.914
$ws.Range("I102").Value = 54118.26
$ws.Range("K102").Value = 54118.26
$ws.Range("M102").Value = -52496.26

# ARM row 122
$ws.Range("H122").Value = 2397.8628
$ws.Range("I122").Value = 2498.6943
$ws.Range("J122").Value = 2155.8667
$ws.Range("K122").Value = 7496.0829
$ws.Range("L122").Value = 6467.6001
$ws.Range("M122").Value = -5046.0829
$ws.Range("N122").Value = -11367.6001

# ARM row 136
$ws.Range("H136").Value = 1142.2609
$ws.Range("I136").Value = 739.0454999999999
$ws.Range("J136").Value = 10013
$ws.Range("K136").Value = 2217.1365
$ws.Range("L136").Value = 30039
$ws.Range("M136").Value = 332.8635000000004
$ws.Range("N136").Value = -35139

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86
$ws.Range("H86").Value = 1530.4667
$ws.Range("I86").Value = 1446
$ws.Range("J86").Value = 1699.4
$ws.Range("K86").Value = 1446
$ws.Range("L86").Value = 1699.4
$ws.Range("M86").Value = -323
$ws.Range("N86").Value = -3945.4

# BSM row 89
$ws.Range("H89").Value = 1530.4667
$ws.Range("I89").Value = 1446
$ws.Range("J89").Value = 1699.4
$ws.Range("K89").Value = 7230
$ws.Range("L89").Value = 8497
$ws.Range("M89").Value = -1614
$ws.Range("N89").Value = -19729

$ws = $wb.Worksheets.Item("CRP")
# CRP row 99
$ws.Range("H99").Value = 1205438.2
$ws.Range("I99").Value = 3188.4614
$ws.Range("J99").Value = 2407688.2
$ws.Range("K99").Value = 3188.4614
$ws.Range("L99").Value = 2407688.2
$ws.Range("M99").Value = -1690.4614
$ws.Range("N99").Value = -2410684.2

# CRP row 105
$ws.Range("H105").Value = 2468.2
$ws.Range("I105").Value = 1743.4375
$ws.Range("J105").Value = 3756.6667
$ws.Range("K105").Value = 1743.4375
$ws.Range("L105").Value = 3756.6667
$ws.Range("M105").Value = 3.5625
$ws.Range("N105").Value = -7250.6667

# CRP row 126
$ws.Range("H126").Value = 1205438.2
$ws.Range("I126").Value = 3188.4614
$ws.Range("J126").Value = 2407688.2
$ws.Range("K126").Value = 9565.3842
$ws.Range("L126").Value = 7223064.600000001
$ws.Range("M126").Value = -7095.3842
$ws.Range("N126").Value = -7228004.600000001

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122
$ws.Range("H122").Value = 10550.714
$ws.Range("I122").Value = 21281.2
$ws.Range("K122").Value = 63843.60000000001
$ws.Range("M122").Value = -61393.60000000001

$ws = $wb.Worksheets.Item("LTW")
# LTW row 132
$ws.Range("H132").Value = 133300
$ws.Range("I132").Value = 133300
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 399900
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -397370
$ws.Range("N132").ClearContents()

# LTW row 136
$ws.Range("H136").Value = 2736.9138
$ws.Range("I136").Value = 2268.353
$ws.Range("J136").Value = 3400.7083
$ws.Range("K136").Value = 6805.059
$ws.Range("L136").Value = 10202.1249
$ws.Range("M136").Value = -4255.059
$ws.Range("N136").Value = -15302.1249

$ws = $wb.Worksheets.Item("WVR")
# WVR row 126
$ws.Range("H126").Value = 2158.6453
$ws.Range("I126").Value = 1806.5
$ws.Range("J126").Value = 2798.9092
$ws.Range("K126").Value = 5419.5
$ws.Range("L126").Value = 8396.7276
$ws.Range("M126").Value = -2949.5
$ws.Range("N126").Value = -13336.7276

# WVR row 136
$ws.Range("H136").Value = 1527.4182
$ws.Range("I136").Value = 1429.2667
$ws.Range("J136").Value = 1969.1
$ws.Range("K136").Value = 4287.800099999999
$ws.Range("L136").Value = 5907.299999999999
$ws.Range("M136").Value = -1737.800099999999
$ws.Range("N136").Value = -11007.3
